$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy formatting (style s="1", row height, etc.) from the last populated
# row (18) down into the new rows 19-23, matching the existing table look.
$ws.Range("A18:B18").Copy()
$ws.Range("A19:B23").PasteSpecial(-4122)

# Row 19: TogglePanel / 显示/隐藏列表
$ws.Range("A19").Value = "TogglePanel"
$ws.Range("B19").Value = "显示/隐藏列表"

# Row 20: common / 中文
$ws.Range("A20").Value = "common"
$ws.Range("B20").Value = "中文"

# Row 21: en / 英语
$ws.Range("A21").Value = "en"
$ws.Range("B21").Value = "英语"

# Row 22: de / 德语
$ws.Range("A22").Value = "de"
$ws.Range("B22").Value = "德语"

# Row 23: zh-CN / 中文
$ws.Range("A23").Value = "zh-CN"
$ws.Range("B23").Value = "中文"
